$wb = $excel.ActiveWorkbook

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 1355
$ws.Cells.Item(92, 9).Value = 1281.25
$ws.Cells.Item(92, 10).Value = 1650
$ws.Cells.Item(92, 11).Value = 1281.25
$ws.Cells.Item(92, 12).Value = 1650
$ws.Cells.Item(92, 13).Value = -33.25
$ws.Cells.Item(92, 14).Value = -4146

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 713.125
$ws.Cells.Item(107, 9).Value = 699.9167
$ws.Cells.Item(107, 11).Value = 699.9167
$ws.Cells.Item(107, 13).Value = 1220.0833

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 7567.4
$ws.Cells.Item(132, 9).Value = 7519.3335
$ws.Cells.Item(132, 10).Value = 8000
$ws.Cells.Item(132, 11).Value = 22558.0005
$ws.Cells.Item(132, 12).Value = 24000
$ws.Cells.Item(132, 13).Value = -20028.0005
$ws.Cells.Item(132, 14).Value = -29060

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 45269.39
$ws.Cells.Item(137, 9).Value = 1614.7693
$ws.Cells.Item(137, 10).Value = 102020.4
$ws.Cells.Item(137, 11).Value = 4844.3079
$ws.Cells.Item(137, 12).Value = 306061.2
$ws.Cells.Item(137, 13).Value = -2294.3079
$ws.Cells.Item(137, 14).Value = -311161.2

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2033.4736
$ws.Cells.Item(61, 9).Value = 1264.5714
$ws.Cells.Item(61, 10).Value = 4186.4
$ws.Cells.Item(61, 11).Value = 1264.5714
$ws.Cells.Item(61, 12).Value = 4186.4
$ws.Cells.Item(61, 13).Value = -1052.5714
$ws.Cells.Item(61, 14).Value = -4610.4

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 4251.8335
$ws.Cells.Item(74, 9).Value = 4337
$ws.Cells.Item(74, 11).Value = 4337
$ws.Cells.Item(74, 13).Value = -3463

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 4251.8335
$ws.Cells.Item(77, 9).Value = 4337
$ws.Cells.Item(77, 11).Value = 21685
$ws.Cells.Item(77, 13).Value = -17317

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 18111.129
$ws.Cells.Item(132, 9).Value = 1872.2084
$ws.Cells.Item(132, 10).Value = 73787.42999999999
$ws.Cells.Item(132, 11).Value = 5616.6252
$ws.Cells.Item(132, 12).Value = 221362.29
$ws.Cells.Item(132, 13).Value = -3086.6252
$ws.Cells.Item(132, 14).Value = -226422.29

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 2033.4736
$ws.Cells.Item(136, 9).Value = 1264.5714
$ws.Cells.Item(136, 10).Value = 4186.4
$ws.Cells.Item(136, 11).Value = 3793.7142
$ws.Cells.Item(136, 12).Value = 12559.2
$ws.Cells.Item(136, 13).Value = -1243.7142
$ws.Cells.Item(136, 14).Value = -17659.2

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 26121.166
$ws.Cells.Item(134, 10).Value = 1330.5714
$ws.Cells.Item(134, 12).Value = 3991.7142
$ws.Cells.Item(134, 14).Value = -9061.7142

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 18510.5
$ws.Cells.Item(31, 10).Value = 6500
$ws.Cells.Item(31, 12).Value = 6500
$ws.Cells.Item(31, 14).Value = -7090

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 18510.5
$ws.Cells.Item(34, 10).Value = 6500
$ws.Cells.Item(34, 12).Value = 6500
$ws.Cells.Item(34, 14).Value = -6904

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 29027.223
$ws.Cells.Item(58, 9).Value = 1374.375
$ws.Cells.Item(58, 10).Value = 250250
$ws.Cells.Item(58, 11).Value = 1374.375
$ws.Cells.Item(58, 12).Value = 250250
$ws.Cells.Item(58, 13).Value = -1171.375
$ws.Cells.Item(58, 14).Value = -250656

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1503.5
$ws.Cells.Item(107, 9).Value = 1068.375
$ws.Cells.Item(107, 11).Value = 1068.375
$ws.Cells.Item(107, 13).Value = 851.625

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 33402.35
$ws.Cells.Item(132, 9).Value = 41125.152
$ws.Cells.Item(132, 11).Value = 123375.456
$ws.Cells.Item(132, 13).Value = -120845.456

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1089.0625
$ws.Cells.Item(134, 9).Value = 866.5333000000001
$ws.Cells.Item(134, 11).Value = 2599.5999
$ws.Cells.Item(134, 13).Value = -64.59990000000016

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 29027.223
$ws.Cells.Item(136, 9).Value = 1374.375
$ws.Cells.Item(136, 10).Value = 250250
$ws.Cells.Item(136, 11).Value = 4123.125
$ws.Cells.Item(136, 12).Value = 750750
$ws.Cells.Item(136, 13).Value = -1573.125
$ws.Cells.Item(136, 14).Value = -755850

# CUL row 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 162.85715
$ws.Cells.Item(13, 9).Value = 108
$ws.Cells.Item(13, 10).Value = 300
$ws.Cells.Item(13, 11).Value = 324
$ws.Cells.Item(13, 12).Value = 900
$ws.Cells.Item(13, 13).Value = -156
$ws.Cells.Item(13, 14).Value = -1236

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 17774
$ws.Cells.Item(87, 9).Value = 9752.333000000001
$ws.Cells.Item(87, 10).Value = 27400
$ws.Cells.Item(87, 11).Value = 29256.999
$ws.Cells.Item(87, 12).Value = 82200
$ws.Cells.Item(87, 13).Value = -28008.999
$ws.Cells.Item(87, 14).Value = -84696

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(90, 8).Value = 17774
$ws.Cells.Item(90, 9).Value = 9752.333000000001
$ws.Cells.Item(90, 10).Value = 27400
$ws.Cells.Item(90, 11).Value = 87770.997
$ws.Cells.Item(90, 12).Value = 246600
$ws.Cells.Item(90, 13).Value = -81530.997
$ws.Cells.Item(90, 14).Value = -259080

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 775.71
$ws.Cells.Item(131, 10).Value = 786.1684
$ws.Cells.Item(131, 12).Value = 2358.5052
$ws.Cells.Item(131, 14).Value = -12438.5052

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 21722.5
$ws.Cells.Item(57, 10).Value = 21722.5
$ws.Cells.Item(57, 12).Value = 21722.5
$ws.Cells.Item(57, 14).Value = -23362.5

# LTW row 50
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).ClearContents()

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3508.5
$ws.Cells.Item(61, 9).Value = 1636.25
$ws.Cells.Item(61, 11).Value = 1636.25
$ws.Cells.Item(61, 13).Value = -1434.25

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 3508.5
$ws.Cells.Item(113, 9).Value = 1636.25
$ws.Cells.Item(113, 11).Value = 1636.25
$ws.Cells.Item(113, 13).Value = 533.75

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4057.4285
$ws.Cells.Item(132, 9).Value = 3376
$ws.Cells.Item(132, 10).Value = 4966
$ws.Cells.Item(132, 11).Value = 10128
$ws.Cells.Item(132, 12).Value = 14898
$ws.Cells.Item(132, 13).Value = -7598
$ws.Cells.Item(132, 14).Value = -19958

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 168333.67
$ws.Cells.Item(136, 9).Value = 168333.67
$ws.Cells.Item(136, 11).Value = 505001.01
$ws.Cells.Item(136, 13).Value = -502451.01

# WVR row 5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 8001.5
$ws.Cells.Item(5, 9).Value = 6001
$ws.Cells.Item(5, 11).Value = 6001
$ws.Cells.Item(5, 13).Value = -5889

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1895029.2
$ws.Cells.Item(107, 9).Value = 812.75
$ws.Cells.Item(107, 10).Value = 3789245.8
$ws.Cells.Item(107, 11).Value = 2438.25
$ws.Cells.Item(107, 12).Value = 11367737.4
$ws.Cells.Item(107, 13).Value = -518.25
$ws.Cells.Item(107, 14).Value = -11371577.4

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2666
$ws.Cells.Item(132, 10).Value = 3499
$ws.Cells.Item(132, 12).Value = 10497
$ws.Cells.Item(132, 14).Value = -15557

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 30304738
$ws.Cells.Item(136, 9).Value = 45456150
$ws.Cells.Item(136, 10).Value = 1908.909
$ws.Cells.Item(136, 11).Value = 136368450
$ws.Cells.Item(136, 12).Value = 5726.727000000001
$ws.Cells.Item(136, 13).Value = -136365900
$ws.Cells.Item(136, 14).Value = -10826.727
